$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 2 with the first game server's data.
# Order matters: it drives the shared-string table build order so it
# matches the target file (127.0.0.1, then GameServer_1, then 000104001).
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("A2").Value = "GameServer_1"
$ws.Range("C2").Value = "GameServer_1"
$ws.Range("B2").Value = "000104001"
$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 4001

# C2/F2 are brand new cells - give them the same text format as the rest
# of the row (style index 1 / numFmtId 49, "@") so IDs like 000104001
# keep their leading zeros and match styling of A2/B2.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"

# Now that F2 holds real data, the list-validation should only start
# applying from F3 downward.
$ws.Range("F2:F1048576").Validation.Delete()
$ws.Range("F3:F1048576").Validation.Add(3, 1, 1, "`"TRUE,FALSE`"")

# Move/collapse the active selection onto G3.
$ws.Range("G3").Select() | Out-Null
